$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "리코"
$ws.Range("C17").Value = "Rico"
$ws.Range("D17").Value = "Rico"
$ws.Range("E17").Value = "Rico"
$ws.Range("A18").Value = "EOF"

$ws.Range("E17").Select()
